$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 updates
$ws.Range("G5").Value = 1.57
$ws.Range("H5").Value = 3.8
$ws.Range("I5").Value = 6
$ws.Range("K5").Value = 2.05
$ws.Range("M5").Value = 1.08
$ws.Range("N5").Value = 7.5
$ws.Range("O5").Value = 1.44
$ws.Range("P5").Value = 2.63
$ws.Range("Q5").Value = 2.35
$ws.Range("R5").Value = 1.57
$ws.Range("U5").Value = 2.38
$ws.Range("V5").Value = 1.53
$ws.Range("Z5").Value = 11
$ws.Range("AA5").Value = 17
$ws.Range("AD5").Value = 7.5
$ws.Range("AH5").Value = 11
$ws.Range("AK5").Value = 67
$ws.Range("AO5").Value = 8.5
$ws.Range("AQ5").Value = 29

# Row 10 updates
$ws.Range("O10").Value = 1.36
$ws.Range("P10").Value = 3
$ws.Range("Q10").Value = 2.15
$ws.Range("R10").Value = 1.67

# Row 11 updates
$ws.Range("M11").Value = 1.08
$ws.Range("N11").Value = 8
$ws.Range("Q11").Value = 2.4
$ws.Range("R11").Value = 1.53

# Row 12 updates
$ws.Range("N12").Value = 9
